$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1139.525
$ws.Range("J17").Value = 1167.1333
$ws.Range("L17").Value = 3501.3999
$ws.Range("N17").Value = -3837.3999
$ws.Range("H62").Value = 4298
$ws.Range("I62").Value = 4225
$ws.Range("K62").Value = 4225
$ws.Range("M62").Value = -3601
$ws.Range("H64").Value = 4052.8333
$ws.Range("I64").Value = 4378.5713
$ws.Range("J64").Value = 3596.8
$ws.Range("K64").Value = 4378.5713
$ws.Range("L64").Value = 3596.8
$ws.Range("M64").Value = -4130.5713
$ws.Range("N64").Value = -4092.8
$ws.Range("H65").Value = 4298
$ws.Range("I65").Value = 4225
$ws.Range("K65").Value = 21125
$ws.Range("M65").Value = -18005
$ws.Range("H67").Value = 4052.8333
$ws.Range("I67").Value = 4378.5713
$ws.Range("J67").Value = 3596.8
$ws.Range("K67").Value = 4378.5713
$ws.Range("L67").Value = 3596.8
$ws.Range("M67").Value = -3520.5713
$ws.Range("N67").Value = -5312.8
$ws.Range("H80").Value = 2152.2083
$ws.Range("I80").Value = 1876.5
$ws.Range("J80").Value = 2349.1428
$ws.Range("K80").Value = 5629.5
$ws.Range("L80").Value = 7047.428400000001
$ws.Range("M80").Value = -4631.5
$ws.Range("N80").Value = -9043.428400000001
$ws.Range("H83").Value = 2152.2083
$ws.Range("I83").Value = 1876.5
$ws.Range("J83").Value = 2349.1428
$ws.Range("K83").Value = 16888.5
$ws.Range("L83").Value = 21142.2852
$ws.Range("M83").Value = -11896.5
$ws.Range("N83").Value = -31126.2852
$ws.Range("H86").Value = 16320.692
$ws.Range("I86").Value = 14742.637
$ws.Range("K86").Value = 14742.637
$ws.Range("M86").Value = -13619.637
$ws.Range("H89").Value = 16320.692
$ws.Range("I89").Value = 14742.637
$ws.Range("K89").Value = 73713.185
$ws.Range("M89").Value = -68097.185
$ws.Range("H132").Value = 9742
$ws.Range("I132").Value = 11680.223
$ws.Range("K132").Value = 35040.669
$ws.Range("M132").Value = -32510.669
$ws.Range("H138").Value = 2808.342
$ws.Range("I138").Value = 2261.9312
$ws.Range("K138").Value = 6785.7936
$ws.Range("M138").Value = -1645.7936
$ws.Range("H141").Value = 6399.552
$ws.Range("I141").Value = 2689.85
$ws.Range("J141").Value = 14643.333
$ws.Range("K141").Value = 8069.549999999999
$ws.Range("L141").Value = 43929.999
$ws.Range("M141").Value = -2889.549999999999
$ws.Range("N141").Value = -54289.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 270288.38
$ws.Range("I32").Value = 274890.75
$ws.Range("K32").Value = 274890.75
$ws.Range("M32").Value = -274603.75
$ws.Range("H61").Value = 3134.4
$ws.Range("I61").Value = 2678.2812
$ws.Range("J61").Value = 7999.6665
$ws.Range("K61").Value = 2678.2812
$ws.Range("L61").Value = 7999.6665
$ws.Range("M61").Value = -2466.2812
$ws.Range("N61").Value = -8423.666499999999
$ws.Range("H74").Value = 6515
$ws.Range("I74").Value = 3452.55
$ws.Range("J74").Value = 20125.889
$ws.Range("K74").Value = 3452.55
$ws.Range("L74").Value = 20125.889
$ws.Range("M74").Value = -2578.55
$ws.Range("N74").Value = -21873.889
$ws.Range("H77").Value = 6515
$ws.Range("I77").Value = 3452.55
$ws.Range("J77").Value = 20125.889
$ws.Range("K77").Value = 17262.75
$ws.Range("L77").Value = 100629.445
$ws.Range("M77").Value = -12894.75
$ws.Range("N77").Value = -109365.445
$ws.Range("H122").Value = 29413634
$ws.Range("I122").Value = 45455916
$ws.Range("J122").Value = 2781.3333
$ws.Range("K122").Value = 136367748
$ws.Range("L122").Value = 8343.999899999999
$ws.Range("M122").Value = -136365298
$ws.Range("N122").Value = -13243.9999
$ws.Range("H136").Value = 3134.4
$ws.Range("I136").Value = 2678.2812
$ws.Range("J136").Value = 7999.6665
$ws.Range("K136").Value = 8034.8436
$ws.Range("L136").Value = 23998.9995
$ws.Range("M136").Value = -5484.8436
$ws.Range("N136").Value = -29098.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 80000
$ws.Range("J51").Value = 80000
$ws.Range("L51").Value = 80000
$ws.Range("N51").Value = -80982
$ws.Range("H80").Value = 532.61536
$ws.Range("I80").Value = 725
$ws.Range("J80").Value = 391.53333
$ws.Range("K80").Value = 725
$ws.Range("L80").Value = 391.53333
$ws.Range("M80").Value = 273
$ws.Range("N80").Value = -2387.53333
$ws.Range("H83").Value = 532.61536
$ws.Range("I83").Value = 725
$ws.Range("J83").Value = 391.53333
$ws.Range("K83").Value = 3625
$ws.Range("L83").Value = 1957.66665
$ws.Range("M83").Value = 1367
$ws.Range("N83").Value = -11941.66665
$ws.Range("H105").Value = 19984.5
$ws.Range("J105").Value = 19984.5
$ws.Range("L105").Value = 19984.5
$ws.Range("N105").Value = -23478.5
$ws.Range("H106").Value = 30750
$ws.Range("J106").Value = 30750
$ws.Range("L106").Value = 30750
$ws.Range("N106").Value = -33274
$ws.Range("H134").Value = 2094.9
$ws.Range("I134").Value = 1960.2413
$ws.Range("K134").Value = 5880.7239
$ws.Range("M134").Value = -3345.7239
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7174.684
$ws.Range("I105").Value = 7430.5293
$ws.Range("K105").Value = 7430.5293
$ws.Range("M105").Value = -5683.5293
$ws.Range("H132").Value = 2194.0227
$ws.Range("I132").Value = 2175.279
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6525.837
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3995.837
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 45000
$ws.Range("H121").Value = 16597.176
$ws.Range("I121").Value = 246.75
$ws.Range("K121").Value = 740.25
$ws.Range("M121").Value = 569.75
$ws.Range("H139").Value = 8820.546
$ws.Range("I139").Value = 6139.222
$ws.Range("J139").Value = 10676.846
$ws.Range("K139").Value = 18417.666
$ws.Range("L139").Value = 32030.538
$ws.Range("M139").Value = -13277.666
$ws.Range("N139").Value = -42310.538
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2736.75
$ws.Range("I122").Value = 1998.8
$ws.Range("J122").Value = 3966.6667
$ws.Range("K122").Value = 5996.4
$ws.Range("L122").Value = 11900.0001
$ws.Range("M122").Value = -3546.4
$ws.Range("N122").Value = -16800.0001
$ws.Range("H130").Value = 56487.5
$ws.Range("J130").Value = 56487.5
$ws.Range("L130").Value = 56487.5
$ws.Range("N130").Value = -66527.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3608.0908
$ws.Range("I93").Value = 1343.375
$ws.Range("K93").Value = 1343.375
$ws.Range("M93").Value = -95.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 290475.16
$ws.Range("I4").Value = 337637.75
$ws.Range("J4").Value = 7499.5
$ws.Range("K4").Value = 337637.75
$ws.Range("L4").Value = 7499.5
$ws.Range("M4").Value = -337524.75
$ws.Range("N4").Value = -7725.5
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("H62").Value = 5830.1904
$ws.Range("I62").Value = 5044.091
$ws.Range("J62").Value = 6694.9
$ws.Range("K62").Value = 5044.091
$ws.Range("L62").Value = 6694.9
$ws.Range("M62").Value = -4420.091
$ws.Range("N62").Value = -7942.9
$ws.Range("H65").Value = 5830.1904
$ws.Range("I65").Value = 5044.091
$ws.Range("J65").Value = 6694.9
$ws.Range("K65").Value = 25220.455
$ws.Range("L65").Value = 33474.5
$ws.Range("M65").Value = -22100.455
$ws.Range("N65").Value = -39714.5
$ws.Range("H107").Value = 2265.0571
$ws.Range("I107").Value = 1094.3914
$ws.Range("K107").Value = 3283.1742
$ws.Range("M107").Value = -1363.1742
$ws.Range("H132").Value = 2564.5435
$ws.Range("I132").Value = 1865.2812
$ws.Range("J132").Value = 4162.857
$ws.Range("K132").Value = 5595.8436
$ws.Range("L132").Value = 12488.571
$ws.Range("M132").Value = -3065.8436
$ws.Range("N132").Value = -17548.571
$ws.Range("H136").Value = 1556.5428
$ws.Range("I136").Value = 1180.25
$ws.Range("K136").Value = 3540.75
$ws.Range("M136").Value = -990.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N132").Value = -14060
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N96").Value = -49118
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M7").ClearContents()
